# Update (Analyze PO & Forecast)
# Applies updated forecast figures to the "Forecast Comparison" sheet
# and the corresponding summary statistics on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# --- Forecast Comparison sheet: MyForecast (column D) updates ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("D2").Value = 54
$wsForecast.Range("D3").Value = 45
$wsForecast.Range("D4").Value = 38
$wsForecast.Range("D5").Value = 39
$wsForecast.Range("D6").Value = 52
$wsForecast.Range("D16").Value = 50
$wsForecast.Range("D17").Value = 45

# --- Summary sheet updates ---
# These cells store numeric-/date-looking values as plain TEXT in the
# original file. Force the number format to Text first so Excel's
# auto-type-detection doesn't turn them into real numbers/dates.
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9:B15").NumberFormat = "@"

$wsSummary.Range("B9").Value = "792"
$wsSummary.Range("B10").Value = "389"
$wsSummary.Range("B11").Value = "176"
$wsSummary.Range("B12").Value = "55"
$wsSummary.Range("B13").Value = "2025-03-16"
$wsSummary.Range("B14").Value = "38"
$wsSummary.Range("B15").Value = "2025-02-09"
